# Fruta / hortaliza, semanal
# A new weekly price record was added to the "Cebollín" (Vega Modelo de
# Temuco) data set. In the canonical OOXML this shows up as a brand new
# row 183 with the rest of the existing rows (old 183-217) shifted down
# by one (to 184-218).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 183, shifting rows 183:217 down to 184:218.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(183, 1).Value  = 10
$ws.Cells.Item(183, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(183, 3).Value  = "La Araucanía"
$ws.Cells.Item(183, 4).Value  = 44476
$ws.Cells.Item(183, 5).Value  = 9
$ws.Cells.Item(183, 6).Value  = 100112037
$ws.Cells.Item(183, 7).Value  = "Cebollín"
$ws.Cells.Item(183, 8).Value  = "Sin especificar"
$ws.Cells.Item(183, 9).Value  = "Primera"
$ws.Cells.Item(183, 10).Value = 40
$ws.Cells.Item(183, 11).Value = 7000
$ws.Cells.Item(183, 12).Value = 7000
$ws.Cells.Item(183, 13).Value = 7000
$ws.Cells.Item(183, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(183, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(183, 16).Value = 583
$ws.Cells.Item(183, 17).Value = 12
$ws.Cells.Item(183, 18).Value = "Hortaliza"
